$d = $word.ActiveDocument

# The embedded picture's alt-text (serialized as the "descr" attribute on
# both wp:docPr and pic:cNvPr in the underlying OOXML) holds a generated
# temp file name of the form ".../anydsl%20class%20diagram_copy-m2doc<digits>.jpg".
# Update it to reflect the newly generated temp file name produced by the
# refreshed image export.
$oldDescr = "file:/tmp/anydsl%20class%20diagram_copy-m2doc9004668813486921900.jpg"
$newDescr = "file:/tmp/anydsl%20class%20diagram_copy-m2doc1427643822779169508.jpg"

foreach ($shape in $d.InlineShapes) {
    if ($shape.AlternativeText -eq $oldDescr) {
        $shape.AlternativeText = $newDescr
    }
}
